$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is new; copy formatting from row 18 (column A date style) before writing values
$ws.Range("A18").Copy($ws.Range("A19"))

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 3.145939949069287
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 2.51031180018495

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 1.769627576887389
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = -1.56363396419209

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -4.774178217057779
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 1.68591416918662

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 1.97975191822708
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 3.609042024648068

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 3.452886745653183
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 2.11424984565185

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 1.239479831392853
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 0.9276342348636168

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.2379616621361214
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 1.3307042289459

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 1.51977456621637
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 0.3626364251072101

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 1.470039379455756
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 1.339087911421144

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 1.638797242243251
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 1.006353890555189

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 2.161565493242668
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 3.257358596620663

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 2.214251681313772
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = -0.3013396321239648

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 0.6066442151010376
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 0.2691345740889695

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -4.207901339433196
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = -0.6332519459683494

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 1.099928004397532
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 0.4297312830033428

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 2.310042359896225
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = -0.177017417229286

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = 0.0464415346324687
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.08221134935635366

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -0.3101476031197148
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 0.5295174046934692
